$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(3,4,5,7,8,10,11,12,13,14,15,17,18,19,21,22,23,24,25,26,27,28,30,32,34,35,36,39,41,42,44,47)

foreach ($r in $rows) {
    $ws.Range("B$r`:J$r").Value = 0
}
